$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: rate ratio value in C2 should be 2.3 (was mistakenly 2.29)
$ws.Range("C2").Value = 2.3

# Fix typo: value in D5 should be 5.3 (was mistakenly 53)
$ws.Range("D5").Value = 5.3

# Update the active selection to C3
$ws.Range("C3").Select()
